$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4331.6665
$ws.Range("I32").Value = 3999
$ws.Range("J32").Value = 4498
$ws.Range("K32").Value = 3999
$ws.Range("L32").Value = 4498
$ws.Range("M32").Value = -3673
$ws.Range("N32").Value = -5150
$ws.Range("H43").Value = 5857
$ws.Range("I43").Value = 9000
$ws.Range("K43").Value = 9000
$ws.Range("M43").Value = -8931
$ws.Range("H53").Value = 437.54544
$ws.Range("I53").Value = 208.85715
$ws.Range("K53").Value = 208.85715
$ws.Range("M53").Value = 428.14285
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H137").Value = 2512.2104
$ws.Range("I137").Value = 2159.3125
$ws.Range("K137").Value = 6477.9375
$ws.Range("M137").Value = -3927.9375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5081.0977
$ws.Range("I32").Value = 3062.4866
$ws.Range("K32").Value = 3062.4866
$ws.Range("M32").Value = -2775.4866
$ws.Range("H61").Value = 90910890
$ws.Range("I61").Value = 100001784
$ws.Range("K61").Value = 100001784
$ws.Range("M61").Value = -100001572
$ws.Range("H97").Value = 986.8333
$ws.Range("I97").Value = 984.2
$ws.Range("K97").Value = 984.2
$ws.Range("M97").Value = -488.2
$ws.Range("H132").Value = 2943464
$ws.Range("J132").Value = 3450
$ws.Range("L132").Value = 10350
$ws.Range("N132").Value = -15410
$ws.Range("H136").Value = 90910890
$ws.Range("I136").Value = 100001784
$ws.Range("K136").Value = 300005352
$ws.Range("M136").Value = -300002802
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 450
$ws.Range("I12").Value = 450
$ws.Range("K12").Value = 450
$ws.Range("M12").Value = -282
$ws.Range("H14").Value = 5999.5
$ws.Range("J14").Value = 5999.5
$ws.Range("L14").Value = 5999.5
$ws.Range("N14").Value = -6343.5
$ws.Range("H16").Value = 10009
$ws.Range("J16").Value = 10009
$ws.Range("L16").Value = 10009
$ws.Range("N16").Value = -10349
$ws.Range("H36").Value = 4839.8335
$ws.Range("I36").Value = 4599.6
$ws.Range("J36").Value = 6041
$ws.Range("K36").Value = 4599.6
$ws.Range("L36").Value = 6041
$ws.Range("M36").Value = -4065.6
$ws.Range("N36").Value = -7109
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H52").Value = 70000
$ws.Range("J52").Value = 70000
$ws.Range("L52").Value = 70000
$ws.Range("N52").Value = -70526
$ws.Range("H62").Value = 82900
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 82900
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H105").Value = 2790.1428
$ws.Range("I105").Value = 2145.25
$ws.Range("J105").Value = 3650
$ws.Range("K105").Value = 2145.25
$ws.Range("L105").Value = 3650
$ws.Range("M105").Value = -398.25
$ws.Range("N105").Value = -7144
$ws.Range("H110").Value = 22222
$ws.Range("J110").Value = 22222
$ws.Range("L110").Value = 22222
$ws.Range("N110").Value = -30402
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 70000
$ws.Range("J121").Value = 70000
$ws.Range("L121").Value = 70000
$ws.Range("N121").Value = -73494
$ws.Range("H124").Value = 93412
$ws.Range("J124").Value = 93412
$ws.Range("L124").Value = 93412
$ws.Range("N124").Value = -103232
$ws.Range("H125").Value = 78999
$ws.Range("J125").Value = 78999
$ws.Range("L125").Value = 78999
$ws.Range("N125").Value = -88839
$ws.Range("H129").Value = 198938
$ws.Range("J129").Value = 198938
$ws.Range("L129").Value = 198938
$ws.Range("N129").Value = -208938
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 83352090
$ws.Range("I58").Value = 100021496
$ws.Range("K58").Value = 100021496
$ws.Range("M58").Value = -100021293
$ws.Range("H134").Value = 9300528
$ws.Range("I134").Value = 10462368
$ws.Range("K134").Value = 31387104
$ws.Range("M134").Value = -31384569
$ws.Range("H136").Value = 83352090
$ws.Range("I136").Value = 100021496
$ws.Range("K136").Value = 300064488
$ws.Range("M136").Value = -300061938
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2848.75
$ws.Range("I68").Value = 475
$ws.Range("J68").Value = 5222.5
$ws.Range("K68").Value = 1425
$ws.Range("L68").Value = 15667.5
$ws.Range("M68").Value = -614
$ws.Range("N68").Value = -17289.5
$ws.Range("H71").Value = 2848.75
$ws.Range("I71").Value = 475
$ws.Range("J71").Value = 5222.5
$ws.Range("K71").Value = 4275
$ws.Range("L71").Value = 47002.5
$ws.Range("M71").Value = -219
$ws.Range("N71").Value = -55114.5
$ws.Range("H114").Value = 252974.75
$ws.Range("I114").Value = 999999
$ws.Range("J114").Value = 3966.6667
$ws.Range("K114").Value = 2999997
$ws.Range("L114").Value = 11900.0001
$ws.Range("M114").Value = -2996743
$ws.Range("N114").Value = -18408.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2246.6667
$ws.Range("I40").Value = 2246.6667
$ws.Range("K40").Value = 2246.6667
$ws.Range("M40").Value = -2110.6667
$ws.Range("H61").Value = 5193.5
$ws.Range("I61").Value = 5193.5
$ws.Range("K61").Value = 5193.5
$ws.Range("M61").Value = -4991.5
$ws.Range("H113").Value = 5193.5
$ws.Range("I113").Value = 5193.5
$ws.Range("K113").Value = 5193.5
$ws.Range("M113").Value = -3023.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 42523.5
$ws.Range("I55").Value = 40048
$ws.Range("K55").Value = 40048
$ws.Range("M55").Value = -39771
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H136").Value = 25001858
$ws.Range("I136").Value = 25001858
$ws.Range("K136").Value = 75005574
$ws.Range("M136").Value = -75003024
